$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new data row above the current row 29 (shifts old rows 29-40 down
# to 30-41, which matches every row below unchanged but shifted by one).
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record.
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44489
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 100112026
$ws.Range("G29").Value = "Haba"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 7000
$ws.Range("M29").Value = 7000
$ws.Range("N29").Value = "$/saco 25 kilos"
$ws.Range("O29").Value = "Región de O'Higgins"
$ws.Range("P29").Value = 280
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
